$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# RealTimePlayerData sheet rework ("change RealTimePlayerSata data type"):
# the devil stat table gained three new leading columns (DevilName / Attack
# / HP), the remaining columns shifted over, SoulNumber/ChooseDevil were
# dropped, PrefabPath moved to the end (now column J), Exp became
# ExpEffect, and the whole data row switched from placeholder int/string
# sample values to a real "Reaper" stat line (the old int column - now
# DamageCut - is typed float like everything else).
# -------------------------------------------------------------------------

# Row 1 - English field names
$ws.Range("A1").Value = "DevilName"
$ws.Range("B1").Value = "Attack"
$ws.Range("C1").Value = "HP"
$ws.Range("D1").Value = "Speed"
$ws.Range("E1").Value = "ExpEffect"
$ws.Range("F1").Value = "AbsorbExpRange"
$ws.Range("G1").Value = "DamageCut"
$ws.Range("H1").Value = "Recovery"
$ws.Range("I1").Value = "DropRate"
$ws.Range("J1").Value = "PrefabPath"

# Row 2 - Chinese field names
$ws.Range("A2").Value = "惡魔名稱"
$ws.Range("B2").Value = "攻擊力"
$ws.Range("C2").Value = "血量"
$ws.Range("D2").Value = "速度"
$ws.Range("E2").Value = "經驗值成長效率"
$ws.Range("F2").Value = "吸收經驗範圍"
$ws.Range("G2").Value = "減傷"
$ws.Range("H2").Value = "恢復能力"
$ws.Range("I2").Value = "寶箱掉落率"
$ws.Range("J2").Value = "下載魔王的prefab"

# Row 3 - field types
$ws.Range("A3").Value = "string"
$ws.Range("B3").Value = "float"
$ws.Range("C3").Value = "float"
$ws.Range("D3").Value = "float"
$ws.Range("E3").Value = "float"
$ws.Range("F3").Value = "float"
$ws.Range("G3").Value = "float"
$ws.Range("H3").Value = "float"
$ws.Range("I3").Value = "float"
$ws.Range("J3").Value = "string"

# Row 4 - sample data (Reaper)
$ws.Range("A4").Value = "Reaper"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "Assets/Prefabs/Devils/Reaper.prefab"

# The whole table is left aligned (field-name rows, type row, data row).
$ws.Range("A1:J4").HorizontalAlignment = -4131

# Highlight the header cells that call out the columns that moved the most
# (AbsorbExpRange, PrefabPath) with the new orange fill used on the
# refreshed header row.
$ws.Range("F1").Interior.Color = 39423
$ws.Range("J1").Interior.Color = 39423

# White-fill the tail of the type row (F3:I3), matching the refreshed look
# carried over from the previous "int" column now that everything is float.
$ws.Range("F3:I3").Interior.Color = 16777215

# A couple of trailing helper cells (K column on every row, plus a lone
# marker in column O) are present-but-blank in the refreshed sheet.
$ws.Range("K1:K4").HorizontalAlignment = -4131
$ws.Range("O1:O4").HorizontalAlignment = -4131

# New column widths for I, J, L, M (character units; the sheet uses Arial
# 10 metrics, so inputs are chosen to land on the requested stored width).
$ws.Columns.Item(9).ColumnWidth = 15.1666666667
$ws.Columns.Item(10).ColumnWidth = 32.5
$ws.Columns.Item(12).ColumnWidth = 16.3333333333
$ws.Columns.Item(13).ColumnWidth = 33.6666666667
